$p = $ppt.ActivePresentation

# --- 1. Refresh the cached "datetimeFigureOut" date field text everywhere it
#        appears (slide master + every slide layout's Date placeholder). ---
$newDate = "12/19/2021"
$m = $p.SlideMaster

for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.PlaceholderFormat.Type -eq 16) {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($j = 1; $j -le $m.CustomLayouts.Count; $j++) {
    $layout = $m.CustomLayouts.Item($j)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 2. Rename CheKiPEUQ -> PEUQSE in the closing paragraph of TextBox 22 on
#        slide 1, splitting it out into its own run like the authored edit. ---
$s = $p.Slides.Item(1)
$box = $s.Shapes.Item("TextBox 22")
$tr = $box.TextFrame.TextRange
$full = $tr.Text
$idx = $full.IndexOf("CheKiPEUQ ")
$hit = $tr.Characters($idx + 1, 10)
$hit.Text = "PEUQSE "
